# Update Sprint Burndown Chart
# Fill in the actual "Storypoints Ist" (done-so-far) values for Sprint 1
# on the burndown tracker; the "remaining" (verbleibend) formulas in
# columns C and E recalculate automatically from these inputs, and the
# burndown chart (which plots those columns) follows along.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vokabeltrainer Sprint 1")

$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 10

# Leave the selection on the last entered cell, like a user would.
$ws.Range("B7").Select() | Out-Null
